$wb = $excel.ActiveWorkbook
$ws = $wb.ActiveSheet

# Add a new row for the "Town of Belmere" location with its quest,
# matching the existing table's layout (Location in column A, Quests in column D).
$ws.Range("A12").Value = "Town of Belmere"
$ws.Range("D12").Value = "001012 Burden of Beast"

# Leave the selection where Excel would land after typing into D12 and pressing Enter.
$ws.Range("D13").Select()
